# Update "想去人数" (interested-count) figures in the F column on the
# "展览" and "全部类型" worksheets, per the latest scrape run.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 4445
$ws1.Range("F3").Value  = 2470
$ws1.Range("F5").Value  = 26
$ws1.Range("F8").Value  = 222
$ws1.Range("F10").Value = 163
$ws1.Range("F12").Value = 1663
$ws1.Range("F14").Value = 3589
$ws1.Range("F15").Value = 9
$ws1.Range("F16").Value = 240

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 4445
$ws4.Range("F3").Value  = 2470
$ws4.Range("F5").Value  = 26
$ws4.Range("F10").Value = 222
$ws4.Range("F12").Value = 163
$ws4.Range("F16").Value = 1663
$ws4.Range("F18").Value = 3589
$ws4.Range("F19").Value = 9
$ws4.Range("F20").Value = 240
